# Generate Report for Handoff
# Inserts two new file records (71c9cb71-... and f28872da-...) into the
# localization-status report: one spliced in between the existing
# 57c0ca54-... and df209093-... rows, the other appended at the end.
# Also refreshes the "Latest Handoff"/handoff-datetime timestamps for the
# already-present rows. Applied identically across the Overview, zh-cn and
# de-de sheets (and their tables/hyperlinks).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "71c9cb71-a061-4618-ad79-facce3a818ba.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-12-13 10:12:26"

$ws.Range("A5").Value = "f28872da-d85c-48b2-be4e-420ac0c6ec83.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "2016-12-13 10:12:26"

$ws.Range("D2").Value = "2016-12-13 10:12:26"
$ws.Range("D4").Value = "2016-12-13 10:12:26"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/57c0ca54-76d2-4668-b60d-04cb3c365f29.md", "", "", "57c0ca54-76d2-4668-b60d-04cb3c365f29.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/71c9cb71-a061-4618-ad79-facce3a818ba.md", "", "", "71c9cb71-a061-4618-ad79-facce3a818ba.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/df209093-f042-40bd-9c1b-0b560ad035ef.md", "", "", "df209093-f042-40bd-9c1b-0b560ad035ef.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/f28872da-d85c-48b2-be4e-420ac0c6ec83.md", "", "", "f28872da-d85c-48b2-be4e-420ac0c6ec83.md")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "71c9cb71-a061-4618-ad79-facce3a818ba.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "71c9cb71-a061-4618-ad79-facce3a818ba.a768a007a2c38200454d7cfcd15382cdbbbe1a22.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-13 10:12:21"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

$ws.Range("A5").Value = "f28872da-d85c-48b2-be4e-420ac0c6ec83.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "f28872da-d85c-48b2-be4e-420ac0c6ec83.a42cf9e9755e5b28af566f6070a064c265c630f5.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-13 10:12:21"
$ws.Range("H5").Value = "0001-01-01 00:00:00"
$ws.Range("I5").Value = "Include"

$ws.Range("E2").Value = "2016-03-13 10:12:21"
$ws.Range("E4").Value = "2016-03-13 10:12:21"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/57c0ca54-76d2-4668-b60d-04cb3c365f29.md", "", "", "57c0ca54-76d2-4668-b60d-04cb3c365f29.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/57c0ca54-76d2-4668-b60d-04cb3c365f29.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/572aa069aec0b8af9eb2441422f1413aeb9a04eb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/57c0ca54-76d2-4668-b60d-04cb3c365f29.1266a7c9cb50fd23985b815f424a23d147db8b3f.zh-cn.xlf", "", "", "57c0ca54-76d2-4668-b60d-04cb3c365f29.1266a7c9cb50fd23985b815f424a23d147db8b3f.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/71c9cb71-a061-4618-ad79-facce3a818ba.md", "", "", "71c9cb71-a061-4618-ad79-facce3a818ba.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/71c9cb71-a061-4618-ad79-facce3a818ba.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/572aa069aec0b8af9eb2441422f1413aeb9a04eb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/71c9cb71-a061-4618-ad79-facce3a818ba.a768a007a2c38200454d7cfcd15382cdbbbe1a22.zh-cn.xlf", "", "", "71c9cb71-a061-4618-ad79-facce3a818ba.a768a007a2c38200454d7cfcd15382cdbbbe1a22.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/df209093-f042-40bd-9c1b-0b560ad035ef.md", "", "", "df209093-f042-40bd-9c1b-0b560ad035ef.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/df209093-f042-40bd-9c1b-0b560ad035ef.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/572aa069aec0b8af9eb2441422f1413aeb9a04eb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/df209093-f042-40bd-9c1b-0b560ad035ef.580a67526ec99aae32e02b1887bcbfa22d92dbb8.zh-cn.xlf", "", "", "df209093-f042-40bd-9c1b-0b560ad035ef.580a67526ec99aae32e02b1887bcbfa22d92dbb8.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/f28872da-d85c-48b2-be4e-420ac0c6ec83.md", "", "", "f28872da-d85c-48b2-be4e-420ac0c6ec83.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/f28872da-d85c-48b2-be4e-420ac0c6ec83.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/572aa069aec0b8af9eb2441422f1413aeb9a04eb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/f28872da-d85c-48b2-be4e-420ac0c6ec83.a42cf9e9755e5b28af566f6070a064c265c630f5.zh-cn.xlf", "", "", "f28872da-d85c-48b2-be4e-420ac0c6ec83.a42cf9e9755e5b28af566f6070a064c265c630f5.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "71c9cb71-a061-4618-ad79-facce3a818ba.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "71c9cb71-a061-4618-ad79-facce3a818ba.a768a007a2c38200454d7cfcd15382cdbbbe1a22.de-de.xlf"
$ws.Range("E3").Value = "2016-03-13 10:12:26"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

$ws.Range("A5").Value = "f28872da-d85c-48b2-be4e-420ac0c6ec83.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "f28872da-d85c-48b2-be4e-420ac0c6ec83.a42cf9e9755e5b28af566f6070a064c265c630f5.de-de.xlf"
$ws.Range("E5").Value = "2016-03-13 10:12:26"
$ws.Range("H5").Value = "0001-01-01 00:00:00"
$ws.Range("I5").Value = "Include"

$ws.Range("E2").Value = "2016-03-13 10:12:26"
$ws.Range("E4").Value = "2016-03-13 10:12:26"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/57c0ca54-76d2-4668-b60d-04cb3c365f29.md", "", "", "57c0ca54-76d2-4668-b60d-04cb3c365f29.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/57c0ca54-76d2-4668-b60d-04cb3c365f29.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4bc5e6c9a7bdc9049efbc66a13294004d6cdc708/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/57c0ca54-76d2-4668-b60d-04cb3c365f29.1266a7c9cb50fd23985b815f424a23d147db8b3f.de-de.xlf", "", "", "57c0ca54-76d2-4668-b60d-04cb3c365f29.1266a7c9cb50fd23985b815f424a23d147db8b3f.de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/71c9cb71-a061-4618-ad79-facce3a818ba.md", "", "", "71c9cb71-a061-4618-ad79-facce3a818ba.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/71c9cb71-a061-4618-ad79-facce3a818ba.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4bc5e6c9a7bdc9049efbc66a13294004d6cdc708/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/71c9cb71-a061-4618-ad79-facce3a818ba.a768a007a2c38200454d7cfcd15382cdbbbe1a22.de-de.xlf", "", "", "71c9cb71-a061-4618-ad79-facce3a818ba.a768a007a2c38200454d7cfcd15382cdbbbe1a22.de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/df209093-f042-40bd-9c1b-0b560ad035ef.md", "", "", "df209093-f042-40bd-9c1b-0b560ad035ef.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/df209093-f042-40bd-9c1b-0b560ad035ef.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4bc5e6c9a7bdc9049efbc66a13294004d6cdc708/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/df209093-f042-40bd-9c1b-0b560ad035ef.580a67526ec99aae32e02b1887bcbfa22d92dbb8.de-de.xlf", "", "", "df209093-f042-40bd-9c1b-0b560ad035ef.580a67526ec99aae32e02b1887bcbfa22d92dbb8.de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/f28872da-d85c-48b2-be4e-420ac0c6ec83.md", "", "", "f28872da-d85c-48b2-be4e-420ac0c6ec83.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/18a195d96adffc3ddccc6d8a4839582a0b70a256/e2e/f28872da-d85c-48b2-be4e-420ac0c6ec83.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4bc5e6c9a7bdc9049efbc66a13294004d6cdc708/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/f28872da-d85c-48b2-be4e-420ac0c6ec83.a42cf9e9755e5b28af566f6070a064c265c630f5.de-de.xlf", "", "", "f28872da-d85c-48b2-be4e-420ac0c6ec83.a42cf9e9755e5b28af566f6070a064c265c630f5.de-de.xlf")

# Leave the Overview sheet active/selected, matching the original workbook.
$wb.Worksheets.Item(1).Activate()
